$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: new entry "10" dated 1/17/2022 for RPA RLOGIC, completed
$ws.Range("A18").Value = 10
$ws.Range("B18").NumberFormat = "m/d/yy"
$ws.Range("B18").Value = 44578
$ws.Range("C18").Value = "RPA RLOGIC"
$ws.Range("D18").Value = "1. Removing rows where minus sign at Other Column in the warranty monthly files has been done during append to the previous month file"
$ws.Range("E18").NumberFormat = "0%"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = "Completed"

# Row 19: second comment line for the same entry, completed
$ws.Range("D19").Value = "2. ESA task has been tested success at Mohan san's system and confirmed  for the ESA today"
$ws.Range("E19").NumberFormat = "0%"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = "Completed"

# Update the active selection to match the author's last cursor position
$ws.Range("D25").Select() | Out-Null
